$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("glycan", "binding_score", "monosaccharides", "motifs", "sasa", "flexibility", "has_multi_node_motifs")

$rows = @(
    @("Neu5Ac(a2-3)Gal(b1-4)GlcNAc", 1.389422566655655, "['Neu5Ac(a2-3)', 'Gal(b1-4)', 'GlcNAc(b1-1)']", "['Sia(a2-3)Gal(b1-4)GlcNAc']", 8.799612018495054, 1.396110483550542, $true),
    @("Neu5Ac(a2-3)Gal(b1-4)GlcNAc(b1-2)Man(a1-3)[Neu5Ac(a2-3)Gal(b1-4)GlcNAc(b1-2)Man(a1-6)]Man(b1-4)GlcNAc(b1-4)GlcNAc", -0.2362585490045806, "['Neu5Ac(a2-3)', 'Gal(b1-4)', 'GlcNAc(b1-2)', 'Neu5Ac(a2-3)', 'Gal(b1-4)', 'GlcNAc(b1-2)']", "['Sia(a2-3)Gal(b1-4)GlcNAc']", 15.14707107057198, 23.79159542523948, $true),
    @("Neu5Ac(a2-3)Gal(b1-4)[Fuc(a1-3)]GlcNAc", -0.1503777538251799, "['Neu5Ac(a2-3)', 'Gal(b1-4)', 'GlcNAc(b1-1)']", "['Sia(a2-3)Gal(b1-4)GlcNAc']", 7.726713635242172, 2.063616819885816, $true),
    @("Neu5Ac(a2-3)Gal(b1-4)[Fuc(a1-3)]GlcNAc(b1-3)Gal", -0.4166485332269986, "['Neu5Ac(a2-3)', 'Gal(b1-4)', 'GlcNAc(b1-3)']", "['Sia(a2-3)Gal(b1-4)GlcNAc']", 7.485843955616165, 2.611071119407292, $true),
    @("Neu5Gc(a2-3)Gal(b1-4)GlcNAc", -0.4157136630359397, "['Neu5Gc(a2-3)', 'Gal(b1-4)', 'GlcNAc(b1-1)']", "['Sia(a2-3)Gal(b1-4)GlcNAc']", 8.943534070975648, 1.454662418287213, $true),
    @("Neu5Gc(a2-3)Gal(b1-4)[Fuc(a1-3)]GlcNAc", -0.3709586205474984, "['Neu5Gc(a2-3)', 'Gal(b1-4)', 'GlcNAc(b1-1)']", "['Sia(a2-3)Gal(b1-4)GlcNAc']", 7.910425322070216, 1.956328000680412, $true)
)

# Header row values
for ($j = 0; $j -lt $headers.Length; $j++) {
    $ws.Cells.Item(1, $j + 1).Value = $headers[$j]
}

# Data row values
for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}

# Build the header/index style once on A1 (bold font, thin box border,
# centered horizontal + top vertical alignment), then propagate it by
# copy/paste-special so every styled cell shares the exact same style
# record instead of each Range-level property assignment minting a new one.
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Borders.LineStyle = 1
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160

$a1.Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$a1.Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)

$excel.CutCopyMode = $false
